$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "repo-scanner"
$ws.Range("B2").Value = "https://github.com/j-chaganti/repo-scanner"
